$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ApproveJournal")

for ($r = 2; $r -le 20; $r++) {
    $ws.Cells.Item($r, 1).Value = "ApproveJournal"
}

$ws.Range("B15").Select()
